$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the reaction ID in A2 from "v37" to "co2in"
$ws.Range("A2").Value = "co2in"

# Move the active selection from D5 to D7
[void]$ws.Range("D7").Select()
